$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.477.22"
$ws.Range("E2").Value = '  -0.92%  '

$ws.Range("D3").Value = "'1.654.68"
$ws.Range("E3").Value = '  -3.00%  '

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = '  +0.54%  '

$ws.Range("D5").Value = "'307.30"
$ws.Range("E5").Value = '  -0.79%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("D7").Value = "'0.3605"
$ws.Range("E7").Value = '  -3.59%  '

$ws.Range("D8").Value = "'47.43"
$ws.Range("E8").Value = '  -3.04%  '

$ws.Range("D9").Value = "'0.3239"
$ws.Range("E9").Value = '  -5.92%  '

$ws.Range("D10").Value = "'1.118"
$ws.Range("E10").Value = '  -5.65%  '

$ws.Range("D11").Value = "'0.06932"
$ws.Range("E11").Value = '  -7.13%  '

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("D13").Value = "'5.870"
$ws.Range("E13").Value = '  -5.99%  '

$ws.Range("D14").Value = "'19.30"
$ws.Range("E14").Value = '  -7.79%  '

$ws.Range("D15").Value = "'1.655.77"
$ws.Range("E15").Value = '  -2.92%  '

$ws.Range("D16").Value = "'6.540"
$ws.Range("E16").Value = '  -5.73%  '

$ws.Range("D17").Value = "'0.00001041"
$ws.Range("E17").Value = '  -7.15%  '

$ws.Range("D18").Value = "'0.06528"
$ws.Range("E18").Value = '  -2.76%  '

$ws.Range("D19").Value = "'0.9991"
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("D20").Value = "'76.10"
$ws.Range("E20").Value = '  -8.95%  '

$ws.Range("D21").Value = "'5.895"
$ws.Range("E21").Value = '  -7.10%  '

$ws.Range("D22").Value = "'15.58"
$ws.Range("E22").Value = '  -9.09%  '

$ws.Range("D23").Value = "'12.51"
$ws.Range("E23").Value = '  -5.28%  '

$ws.Range("D24").Value = "'24.479.61"
$ws.Range("E24").Value = '  -0.69%  '

$ws.Range("D25").Value = "'2.461"
$ws.Range("E25").Value = '  +2.08%  '

$ws.Range("D26").Value = "'2.290"
$ws.Range("E26").Value = '  -17.31%  '

$ws.Range("D27").Value = "'146.63"
$ws.Range("E27").Value = '  -2.41%  '

$ws.Range("D28").Value = "'18.39"
$ws.Range("E28").Value = '  -8.83%  '

$ws.Range("D29").Value = "'1.839.57"
$ws.Range("E29").Value = '  -2.92%  '

$ws.Range("D30").Value = "'1.187"
$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("D31").Value = "'123.64"
$ws.Range("E31").Value = '  -5.86%  '

$ws.Range("D32").Value = "'4.052"
$ws.Range("E32").Value = '  -3.46%  '

$ws.Range("D33").Value = "'5.594"
$ws.Range("E33").Value = '  -17.19%  '

$ws.Range("D34").Value = "'1.698"
$ws.Range("E34").Value = '  -4.47%  '

$ws.Range("D35").Value = "'0.08337"
$ws.Range("E35").Value = '  -5.34%  '

$ws.Range("D36").Value = "'12.29"
$ws.Range("E36").Value = '  -10.58%  '

$ws.Range("D37").Value = "'5.175"
$ws.Range("E37").Value = '  -6.33%  '

$ws.Range("D38").Value = "'0.06024"
$ws.Range("E38").Value = '  -7.78%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'8.197"
$ws.Range("E39").Value = '  -8.34%  '

$ws.Range("D40").Value = "'1.201"
$ws.Range("E40").Value = '  -5.82%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = "'0.2046"
$ws.Range("E41").Value = '  -7.81%  '

$ws.Range("D42").Value = "'0.02180"
$ws.Range("E42").Value = '  -8.25%  '

$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").Value = "'0.5868"
$ws.Range("E44").Value = '  -8.50%  '

$ws.Range("D45").Value = "'3.734"
$ws.Range("E45").Value = '  -2.00%  '

$ws.Range("D46").Value = "'12.61"
$ws.Range("E46").Value = '  -9.01%  '

$ws.Range("D47").Value = "'0.5551"
$ws.Range("E47").Value = '  -8.88%  '

$ws.Range("D48").Value = "'121.90"
$ws.Range("E48").Value = '  -5.94%  '

$ws.Range("D49").Value = "'1.927"
$ws.Range("E49").Value = '  -9.01%  '

$ws.Range("D50").Value = "'0.06889"
$ws.Range("E50").Value = '  -5.15%  '

$ws.Range("D51").Value = "'73.82"
$ws.Range("E51").Value = '  -6.80%  '

Write-Host "Updated cryptos list"